$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 12, pushing existing rows 12..98 down to 13..99,
# and carrying formatting (date style on column D etc.) with it.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new record's data.
$ws.Range("A12").Value = 5
$ws.Range("B12").Value = "Macroferia Regional de Talca"
$ws.Range("C12").Value = "Maule"
$ws.Range("D12").Value = 44490
$ws.Range("E12").Value = 7
$ws.Range("F12").Value = 100112031
$ws.Range("G12").Value = "Poroto verde"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 40000
$ws.Range("L12").Value = 40000
$ws.Range("M12").Value = 40000
$ws.Range("N12").Value = "`$/saco 25 kilos"
$ws.Range("O12").Value = "Región del Maule"
$ws.Range("P12").Value = 1600
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = "Hortaliza"
